# Delaware overview workbook - convert numeric "count" cells to text,
# and add a new "Total" summary row to the County sheet.

function Set-TextValue {
    param(
        $Cell,
        [string]$Text
    )
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Sheet "Overall" ---
$ws = $wb.Worksheets.Item("Overall")
Set-TextValue $ws.Cells.Item(2,1) "514"

# --- Sheet "County" ---
$ws = $wb.Worksheets.Item("County")
Set-TextValue $ws.Cells.Item(2,2) "80"
Set-TextValue $ws.Cells.Item(3,2) "334"
Set-TextValue $ws.Cells.Item(4,2) "100"

# New "Total" row (row 5)
Set-TextValue $ws.Cells.Item(5,1) "Total"
Set-TextValue $ws.Cells.Item(5,2) "514"
Set-TextValue $ws.Cells.Item(5,3) '$1,256,309,593'
Set-TextValue $ws.Cells.Item(5,4) "9.93%"
Set-TextValue $ws.Cells.Item(5,5) "-26.82%"
Set-TextValue $ws.Cells.Item(5,6) "73.54%"

# --- Sheet "Congressional District" ---
$ws = $wb.Worksheets.Item("Congressional District")
Set-TextValue $ws.Cells.Item(2,2) "514"
Set-TextValue $ws.Cells.Item(3,2) "514"

# --- Sheet "Size" ---
$ws = $wb.Worksheets.Item("Size")
Set-TextValue $ws.Cells.Item(2,2) "156"
Set-TextValue $ws.Cells.Item(3,2) "151"
Set-TextValue $ws.Cells.Item(4,2) "91"
Set-TextValue $ws.Cells.Item(5,2) "39"
Set-TextValue $ws.Cells.Item(6,2) "49"
Set-TextValue $ws.Cells.Item(7,2) "28"
Set-TextValue $ws.Cells.Item(8,2) "514"

# --- Sheet "Subsector" ---
$ws = $wb.Worksheets.Item("Subsector")
Set-TextValue $ws.Cells.Item(2,2) "48"
Set-TextValue $ws.Cells.Item(3,2) "61"
Set-TextValue $ws.Cells.Item(4,2) "23"
Set-TextValue $ws.Cells.Item(5,2) "33"
Set-TextValue $ws.Cells.Item(6,2) "3"
Set-TextValue $ws.Cells.Item(7,2) "185"
Set-TextValue $ws.Cells.Item(8,2) "4"
Set-TextValue $ws.Cells.Item(9,2) "34"
Set-TextValue $ws.Cells.Item(10,2) "5"
Set-TextValue $ws.Cells.Item(11,2) "114"
Set-TextValue $ws.Cells.Item(12,2) "4"
Set-TextValue $ws.Cells.Item(13,2) "514"

Write-Host "Delaware overview text edits applied."
